$d = $word.ActiveDocument

# 1. Generated date
$d.Content.Find.Execute("Generated: 2026-02-15", $true, $false, $false, $false, $false, $true, 1, $false, "Generated: 2026-02-22", 2)

# 2. Executive summary totals
$d.Content.Find.Execute("Total Federal Climate Resilience Investment: `$656,509,642 across 17 Tribal Nations", $true, $false, $false, $false, $false, $true, 1, $false, "Total Federal Climate Resilience Investment: `$704,792,362 across 18 Tribal Nations", 2)

$d.Content.Find.Execute("Aggregate Economic Impact: `$1,308,869,356 to `$1,745,159,141", $true, $false, $false, $false, $false, $true, 1, $false, "Aggregate Economic Impact: `$1,394,788,252 to `$1,859,717,669", 2)

$d.Content.Find.Execute("Estimated Jobs Supported: 5,817 to 10,907", $true, $false, $false, $false, $false, $true, 1, $false, "Estimated Jobs Supported: 6,199 to 11,623", 2)

$d.Content.Find.Execute("Coverage Gap: 7 of 24 Tribal Nations", $true, $false, $false, $false, $false, $true, 1, $false, "Coverage Gap: 6 of 24 Tribal Nations", 2)

# 3. Regional award landscape section
$d.Content.Find.Execute("Total Federal Climate Resilience Awards: `$656,509,642", $true, $false, $false, $false, $false, $true, 1, $false, "Total Federal Climate Resilience Awards: `$704,792,362", 2)

$d.Content.Find.Execute("Tribal Nations with Awards: 17 of 24 (71%)", $true, $false, $false, $false, $false, $true, 1, $false, "Tribal Nations with Awards: 18 of 24 (75%)", 2)

$d.Content.Find.Execute("Investment Gap: 7 Tribal Nation(s)", $true, $false, $false, $false, $false, $true, 1, $false, "Investment Gap: 6 Tribal Nation(s)", 2)
